$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.004.56"
$ws.Range("E2").Value = "  -1.65%  "

$ws.Range("D3").Value = "2.683.04"
$ws.Range("E3").Value = "  -2.16%  "

$ws.Range("D5").Value = "'553.62"
$ws.Range("E5").Value = "  -3.33%  "

$ws.Range("D6").Value = "'158.65"
$ws.Range("E6").Value = "  -1.30%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.57%  "

$ws.Range("E9").Value = "  -3.43%  "

$ws.Range("E10").Value = "  -2.13%  "

$ws.Range("D11").Value = "'0.368"
$ws.Range("E11").Value = "  -4.42%  "

$ws.Range("D12").Value = "'5.39"
$ws.Range("E12").Value = "  -7.36%  "

$ws.Range("D13").Value = "3.156.32"
$ws.Range("E13").Value = "  -2.19%  "

$ws.Range("D14").Value = "'26.31"
$ws.Range("E14").Value = "  -2.10%  "

$ws.Range("D15").Value = "62.872.97"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("E16").Value = "  -2.51%  "

$ws.Range("D17").Value = "2.681.70"
$ws.Range("E17").Value = "  -1.50%  "

$ws.Range("D18").Value = "'11.90"
$ws.Range("E18").Value = "  -2.31%  "

$ws.Range("D19").Value = "'4.62"
$ws.Range("E19").Value = "  -3.95%  "

$ws.Range("D20").Value = "'345.26"
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = "  -4.59%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").Value = "'0.508"
$ws.Range("E23").Value = "  -2.77%  "

$ws.Range("D24").Value = "'63.33"
$ws.Range("E24").Value = "  -1.60%  "

$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("D27").Value = "'8.21"
$ws.Range("E27").Value = "  -2.65%  "

$ws.Range("D28").Value = "0.0₃0859"
$ws.Range("E28").Value = "  -6.53%  "

$ws.Range("D29").Value = "'1.37"
$ws.Range("E29").Value = "  +2.26%  "

$ws.Range("D30").Value = "'7.26"
$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("E31").Value = "  -1.43%  "

$ws.Range("D32").Value = "'165.13"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").Value = "'4.89"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").Value = "'1.49"
$ws.Range("E34").Value = "  +0.44%  "

$ws.Range("D36").Value = "'19.54"
$ws.Range("E36").Value = "  -3.03%  "

$ws.Range("E37").Value = "  -1.75%  "

$ws.Range("D38").Value = "'349.16"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'6.35"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").Value = "'0.960"
$ws.Range("E40").Value = "  -3.34%  "

$ws.Range("E41").Value = "  -2.67%  "

$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D43").Value = "'20.36"
$ws.Range("E43").Value = "  -3.93%  "

$ws.Range("D44").Value = "'20.83"
$ws.Range("E44").Value = "  -5.40%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.618"
$ws.Range("E45").Value = "  -1.56%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0562"
$ws.Range("E46").Value = "  -4.06%  "

$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").Value = "'11.04"
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("D49").Value = "'0.0974"
$ws.Range("E49").Value = "  -3.26%  "

$ws.Range("E50").Value = "  -3.63%  "

$ws.Range("D51").Value = "2.102.72"
$ws.Range("E51").Value = "  -1.82%  "
